$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two name cells, removing accented characters (data fix)
$ws.Range("A2").Value = "Ignacio Fernandez Fernandez"
$ws.Range("A3").Value = "Nauce Lopez Gonzalez"

# Move the active selection from C2 to A3
$ws.Range("A3").Select()
